# Attendance sheet update: fill in December 9 (column O) attendance for all
# students, and correct a few other cells (N7, N9, N15: A -> L; K18: A -> P).
#
# The "P" / "A" / "L" cells in this sheet are not just plain text - they carry
# specific (bold + colored) direct formatting on top of a shared fill/border.
# Rather than hand-rolling font colors, reuse already-present, untouched
# template cells of each kind (K3 = P, L3 = A, N14 = L) via Copy/PasteSpecial
# (formats only) so the new cells end up with identical effective formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: this runtime's PowerShell subset does not bind named (-Param value)
# arguments on user functions, so call this positionally:
#   Set-AttendanceCell <TargetRef> <TemplateRef> <Value>
function Set-AttendanceCell {
    param(
        [string]$TargetRef,
        [string]$TemplateRef,
        [string]$Value
    )
    $ws.Range($TemplateRef).Copy() | Out-Null
    $ws.Range($TargetRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($TargetRef).Value = $Value
}

# Template (already-correct, untouched) cells for each attendance code.
$P_TEMPLATE = "K3"    # style for "P" (Present)
$A_TEMPLATE = "L3"    # style for "A" (Absent)
$L_TEMPLATE = "N14"   # style for "L" (Leave)

# --- Column O (December 9) for every student row, rows 3-19 ---
Set-AttendanceCell "O3"  $A_TEMPLATE "A"
Set-AttendanceCell "O4"  $P_TEMPLATE "P"
Set-AttendanceCell "O5"  $P_TEMPLATE "P"
Set-AttendanceCell "O6"  $P_TEMPLATE "P"
Set-AttendanceCell "O7"  $A_TEMPLATE "A"
Set-AttendanceCell "O8"  $P_TEMPLATE "P"
Set-AttendanceCell "O9"  $A_TEMPLATE "A"
Set-AttendanceCell "O10" $A_TEMPLATE "A"
Set-AttendanceCell "O11" $P_TEMPLATE "P"
Set-AttendanceCell "O12" $P_TEMPLATE "P"
Set-AttendanceCell "O13" $P_TEMPLATE "P"
Set-AttendanceCell "O14" $P_TEMPLATE "P"
Set-AttendanceCell "O15" $A_TEMPLATE "A"
Set-AttendanceCell "O16" $P_TEMPLATE "P"
Set-AttendanceCell "O17" $A_TEMPLATE "A"
Set-AttendanceCell "O18" $A_TEMPLATE "A"
Set-AttendanceCell "O19" $P_TEMPLATE "P"

# --- Corrections: column N (December 8), A -> L for rows 7, 9, 15 ---
Set-AttendanceCell "N7"  $L_TEMPLATE "L"
Set-AttendanceCell "N9"  $L_TEMPLATE "L"
Set-AttendanceCell "N15" $L_TEMPLATE "L"

# --- Correction: column K (December 5), A -> P for row 18 ---
Set-AttendanceCell "K18" $P_TEMPLATE "P"

# --- Rebuild the header merges (A1:A2, B1:B2, C1:F1, G1:AK1) ---
# Re-merging them (A1:A2 last) reproduces the resulting mergeCells order.
$ws.Range("A1:A2").UnMerge() | Out-Null
$ws.Range("B1:B2").UnMerge() | Out-Null
$ws.Range("C1:F1").UnMerge() | Out-Null
$ws.Range("G1:AK1").UnMerge() | Out-Null
$ws.Range("B1:B2").Merge() | Out-Null
$ws.Range("C1:F1").Merge() | Out-Null
$ws.Range("G1:AK1").Merge() | Out-Null
$ws.Range("A1:A2").Merge() | Out-Null

Write-Host "Attendance sheet updated."
